$d = $word.ActiveDocument

$d.Content.Find.Execute('2025-03-28 Friday', $true, $false, $false, $false, $false, $true, 1, $false, '2025-03-29 Saturday', 2) | Out-Null
$d.Content.Find.Execute('37×21=777', $true, $false, $false, $false, $false, $true, 1, $false, '56×19=1064', 2) | Out-Null
$d.Content.Find.Execute('71×95=6745', $true, $false, $false, $false, $false, $true, 1, $false, '56×20=1120', 2) | Out-Null
$d.Content.Find.Execute('98×98=9604', $true, $false, $false, $false, $false, $true, 1, $false, '27×19=513', 2) | Out-Null
$d.Content.Find.Execute('89×25=2225', $true, $false, $false, $false, $false, $true, 1, $false, '88×32=2816', 2) | Out-Null
$d.Content.Find.Execute('93×34=3162', $true, $false, $false, $false, $false, $true, 1, $false, '80×30=2400', 2) | Out-Null
$d.Content.Find.Execute('32×76=2432', $true, $false, $false, $false, $false, $true, 1, $false, '21×60=1260', 2) | Out-Null
$d.Content.Find.Execute('27×20=540', $true, $false, $false, $false, $false, $true, 1, $false, '52×17=884', 2) | Out-Null
$d.Content.Find.Execute('57×74=4218', $true, $false, $false, $false, $false, $true, 1, $false, '74×72=5328', 2) | Out-Null
$d.Content.Find.Execute('91×12=1092', $true, $false, $false, $false, $false, $true, 1, $false, '81×94=7614', 2) | Out-Null
$d.Content.Find.Execute('61×19=1159', $true, $false, $false, $false, $false, $true, 1, $false, '50×54=2700', 2) | Out-Null
$d.Content.Find.Execute('60×62=3720', $true, $false, $false, $false, $false, $true, 1, $false, '88×18=1584', 2) | Out-Null
$d.Content.Find.Execute('88×46=4048', $true, $false, $false, $false, $false, $true, 1, $false, '90×72=6480', 2) | Out-Null
$d.Content.Find.Execute('36×78=2808', $true, $false, $false, $false, $false, $true, 1, $false, '67×66=4422', 2) | Out-Null
$d.Content.Find.Execute('73×87=6351', $true, $false, $false, $false, $false, $true, 1, $false, '31×27=837', 2) | Out-Null
$d.Content.Find.Execute('71×94=6674', $true, $false, $false, $false, $false, $true, 1, $false, '74×11=814', 2) | Out-Null
$d.Content.Find.Execute('30×31=930', $true, $false, $false, $false, $false, $true, 1, $false, '82×81=6642', 2) | Out-Null
$d.Content.Find.Execute('49×86=4214', $true, $false, $false, $false, $false, $true, 1, $false, '26×42=1092', 2) | Out-Null
$d.Content.Find.Execute('93×18=1674', $true, $false, $false, $false, $false, $true, 1, $false, '19×85=1615', 2) | Out-Null
$d.Content.Find.Execute('89×34=3026', $true, $false, $false, $false, $false, $true, 1, $false, '79×11=869', 2) | Out-Null
$d.Content.Find.Execute('93×65=6045', $true, $false, $false, $false, $false, $true, 1, $false, '47×26=1222', 2) | Out-Null
$d.Content.Find.Execute('56×67=3752', $true, $false, $false, $false, $false, $true, 1, $false, '25×29=725', 2) | Out-Null
$d.Content.Find.Execute('51×98=4998', $true, $false, $false, $false, $false, $true, 1, $false, '94×30=2820', 2) | Out-Null
$d.Content.Find.Execute('35×91=3185', $true, $false, $false, $false, $false, $true, 1, $false, '30×46=1380', 2) | Out-Null
$d.Content.Find.Execute('75×59=4425', $true, $false, $false, $false, $false, $true, 1, $false, '39×39=1521', 2) | Out-Null
$d.Content.Find.Execute('13×68=884', $true, $false, $false, $false, $false, $true, 1, $false, '62×79=4898', 2) | Out-Null
